# carry forward hourly run
#
# On the "All_Scenarios" sheet, column C ("RunMode") toggles between "Yes"
# and "No" for a block of scenario rows so the carry-forward scenarios now
# run in the hourly batch while the previously-hourly block is switched
# back off:
#   - rows 2-42 and 84-90:  RunMode No  -> Yes
#   - rows 131-140 and 207-247: RunMode Yes -> No
# The active selection / visible window is also moved down to the newly
# enabled block (C2:C90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows switching RunMode from "No" to "Yes"
$rowsToYes = @(2..42) + @(84..90)
foreach ($r in $rowsToYes) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# Rows switching RunMode from "Yes" to "No"
$rowsToNo = @(131..140) + @(207..247)
foreach ($r in $rowsToNo) {
    $ws.Cells.Item($r, 3).Value = "No"
}

# Move the selection / view down to the block that is now enabled.
$ws.Activate()
$ws.Range("C2:C90").Select()
